$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1733333333333333
$ws.Range("C2").Value = 0.6088888888888889
$ws.Range("J2").Value = 0.02222222222222222
$ws.Range("P2").Value = 0.1288888888888889
$ws.Range("S2").Value = 0.06666666666666667

# Row 3
$ws.Range("B3").Value = 0.01379310344827586
$ws.Range("C3").Value = 0.05517241379310345
$ws.Range("J3").Value = 0.02068965517241379
$ws.Range("P3").Value = 0.7172413793103448
$ws.Range("S3").Value = 0.1931034482758621

# Row 4
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.7105263157894737
$ws.Range("S4").Value = 0.2368421052631579

# Row 6
$ws.Range("B6").Value = 0.04854368932038835
$ws.Range("D6").Value = 0.01941747572815534
$ws.Range("F6").Value = 0.05825242718446602
$ws.Range("J6").Value = 0.2233009708737864
$ws.Range("O6").Value = 0.004854368932038835
$ws.Range("Q6").Value = 0.1893203883495146
$ws.Range("R6").Value = 0.07281553398058252
$ws.Range("S6").Value = 0.383495145631068

# Row 7
$ws.Range("B7").Value = 0.07317073170731707
$ws.Range("D7").Value = 0.01829268292682927
$ws.Range("E7").Value = 0.006097560975609756
$ws.Range("F7").Value = 0.06097560975609756
$ws.Range("J7").Value = 0.06707317073170732
$ws.Range("O7").Value = 0.0426829268292683
$ws.Range("Q7").Value = 0.1402439024390244
$ws.Range("R7").Value = 0.06707317073170732
$ws.Range("S7").Value = 0.524390243902439

# Row 8
$ws.Range("B8").Value = 0.07474747474747474
$ws.Range("D8").Value = 0.01212121212121212
$ws.Range("F8").Value = 0.07878787878787878
$ws.Range("J8").Value = 0.105050505050505
$ws.Range("O8").Value = 0.01818181818181818
$ws.Range("Q8").Value = 0.1414141414141414
$ws.Range("R8").Value = 0.101010101010101
$ws.Range("S8").Value = 0.4686868686868687

# Row 9
$ws.Range("B9").Value = 0.06951871657754011
$ws.Range("D9").Value = 0.0106951871657754
$ws.Range("F9").Value = 0.08021390374331551
$ws.Range("J9").Value = 0.09625668449197861
$ws.Range("O9").Value = 0.0106951871657754
$ws.Range("Q9").Value = 0.1550802139037433
$ws.Range("R9").Value = 0.09625668449197861
$ws.Range("S9").Value = 0.481283422459893

# Row 10
$ws.Range("B10").Value = 0.1000918273645546
$ws.Range("D10").Value = 0.02112029384756657
$ws.Range("E10").Value = 0.002754820936639119
$ws.Range("F10").Value = 0.06519742883379247
$ws.Range("J10").Value = 0.09550045913682277
$ws.Range("O10").Value = 0.02295684113865932
$ws.Range("Q10").Value = 0.1992653810835629
$ws.Range("R10").Value = 0.09733700642791551
$ws.Range("S10").Value = 0.3957759412304867

# Row 11
$ws.Range("F11").Value = 0.004065040650406504
$ws.Range("G11").Value = 0.1422764227642276
$ws.Range("J11").Value = 0.07317073170731707
$ws.Range("K11").Value = 0.1910569105691057
$ws.Range("L11").Value = 0.5772357723577236
$ws.Range("S11").Value = 0.01219512195121951

# Row 12
$ws.Range("G12").Value = 0.7083333333333334
$ws.Range("J12").Value = 0.2083333333333333
$ws.Range("K12").Value = 0.01388888888888889
$ws.Range("L12").Value = 0.02777777777777778
$ws.Range("S12").Value = 0.04166666666666666

# Row 13
$ws.Range("G13").Value = 0.7
$ws.Range("J13").Value = 0.22
$ws.Range("S13").Value = 0.08

# Row 15
$ws.Range("F15").Value = 0.01388888888888889
$ws.Range("H15").Value = 0.1712962962962963
$ws.Range("I15").Value = 0.06481481481481481
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.05092592592592592
$ws.Range("M15").Value = 0.004629629629629629
$ws.Range("O15").Value = 0.05092592592592592
$ws.Range("S15").Value = 0.2685185185185185

# Row 16
$ws.Range("F16").Value = 0.006289308176100629
$ws.Range("H16").Value = 0.2075471698113208
$ws.Range("I16").Value = 0.0880503144654088
$ws.Range("J16").Value = 0.3396226415094339
$ws.Range("K16").Value = 0.1257861635220126
$ws.Range("M16").Value = 0.03144654088050314
$ws.Range("O16").Value = 0.0440251572327044
$ws.Range("S16").Value = 0.1572327044025157

# Row 17
$ws.Range("F17").Value = 0.01038961038961039
$ws.Range("H17").Value = 0.1948051948051948
$ws.Range("I17").Value = 0.1142857142857143
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.07012987012987013
$ws.Range("M17").Value = 0.02077922077922078
$ws.Range("O17").Value = 0.07272727272727272
$ws.Range("S17").Value = 0.1168831168831169

# Row 18
$ws.Range("F18").Value = 0.02955665024630542
$ws.Range("H18").Value = 0.2463054187192118
$ws.Range("I18").Value = 0.09359605911330049
$ws.Range("J18").Value = 0.3891625615763547
$ws.Range("K18").Value = 0.09359605911330049
$ws.Range("M18").Value = 0.02463054187192118
$ws.Range("N18").Value = 0.004926108374384237
$ws.Range("O18").Value = 0.02955665024630542
$ws.Range("S18").Value = 0.08866995073891626

# Row 19
$ws.Range("F19").Value = 0.01755786113328013
$ws.Range("H19").Value = 0.2458100558659218
$ws.Range("I19").Value = 0.07741420590582601
$ws.Range("J19").Value = 0.3511572226656026
$ws.Range("K19").Value = 0.09018355945730247
$ws.Range("M19").Value = 0.02873104549082203
$ws.Range("N19").Value = 0.0007980845969672786
$ws.Range("O19").Value = 0.07262569832402235
$ws.Range("S19").Value = 0.1157222665602554
